$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.263.81"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "'1.898.67"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'308.04"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.5208"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").Value = "'0.3776"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "'0.07288"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "'21.22"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.9012"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'0.08174"
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("D13").Value = "'96.72"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "'1.898.45"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'5.297"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'0.000008616"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "'14.56"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'27.288.63"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'5.095"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "'10.71"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "'6.415"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'2.303"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'147.30"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'18.24"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "'1.739"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "'115.62"
$ws.Range("D29").Value = "'4.834"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "'4.911"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "'0.09241"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'0.7970"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "'1.232"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "'3.437"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("D36").Value = "'2.965"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").Value = "'2.596"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'0.5677"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'0.01999"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "'8.987"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'6.573"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "'115.34"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("D44").Value = "'0.1517"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "'0.4884"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").Value = "'10.07"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "'63.76"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  +0.51%  "
